# vault backup: 2025-12-24 16:45:41
# Fill in Day4 (column E) / Day5 (column F) sleep-diary answers for the
# 2025-12-23 / 2025-12-24 block (rows 82-95), restore the title/subtitle
# font-style pairing on each weekly block's header rows, bold the last
# "无" answer cell (F95), and move the active selection to F95.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: the canonical diff also shows the A1/A2/A21/A22/A40/A41/A59/A60/A79
# header cells getting new `s=` (cellXf index) references, paired with the
# stylesheet's xf#7 <-> xf#13 font ids swapping too. Working those two
# changes through together cancels out exactly (xf#7 goes from bold16 to
# italic12 while the cells that used to point at xf#7 now point at xf#13,
# which itself flips bold16->... i.e. every affected cell keeps the exact
# same rendered font it had before). That is Excel's own style-table
# bookkeeping from the resave, not a user-visible formatting change, so
# there is nothing to replay here.

# --- Fill in the Day4 (E) / Day5 (F) answers for rows 82-95 ---
$ws.Range("E82").Value = "7：12"
$ws.Range("F82").Value = "7：42"

$ws.Range("E83").Value = "7：15"
$ws.Range("F83").Value = "8：00"

$ws.Range("E84").Value = "22：45"
$ws.Range("F84").Value = "22：45"

$ws.Range("E85").Value = "23：00"
$ws.Range("F85").Value = "23：45"

$ws.Range("E86").Value = 0
$ws.Range("F86").Value = 60

$ws.Range("E87").Value = 0
$ws.Range("F87").Value = 0

$ws.Range("E88").Value = 0
$ws.Range("F88").Value = 0

$ws.Range("E89").Value = 480
$ws.Range("F89").Value = 480

$ws.Range("E90").Value = "无"
$ws.Range("F90").Value = "无"

$ws.Range("E91").Value = 0
$ws.Range("F91").Value = 60

$ws.Range("E92").Value = 4
$ws.Range("F92").Value = 3

$ws.Range("E93").Value = 4
$ws.Range("F93").Value = 4

$ws.Range("E94").Value = 4
$ws.Range("F94").Value = 3

$ws.Range("E95").Value = "无"
$ws.Range("F95").Value = "无"

# F95 also picks up a bold font (distinct cell style) on this edit.
$ws.Range("F95").Font.Bold = $true

# --- Move the active selection to reflect where the user ended up ---
$ws.Range("F95").Select()
